$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A5").Value = "registration"
$ws.Range("B5").Value = "eis_cpa@yahoo.com"
$ws.Range("C5").Value = "Elena"
$ws.Range("D5").Value = "44 Brandywine cir"
$ws.Range("E5").Value = 7324242995
$ws.Range("F5").Value = "spanish"
$ws.Range("G5").Value = "Programming"
$ws.Range("H5").Value = "United States"
$ws.Range("I5").Value = "United States of America"
$ws.Range("J5").Value = 1968
$ws.Range("K5").Value = 2
$ws.Range("L5").Value = 21
$ws.Range("M5").Value = "Elenaau12"

$ws.Columns.Item(2).ColumnWidth = 23.166666666666668
$ws.Columns.Item(4).ColumnWidth = 18
$ws.Columns.Item(5).ColumnWidth = 10.833333333333334
$ws.Columns.Item(7).ColumnWidth = 15.666666666666666
$ws.Columns.Item(8).ColumnWidth = 13.166666666666666
$ws.Columns.Item(9).ColumnWidth = 21.833333333333332
$ws.Columns.Item(10).ColumnWidth = 5.666666666666667
$ws.Columns.Item(11).ColumnWidth = 4.666666666666667
$ws.Columns.Item(12).ColumnWidth = 4.166666666666667
$ws.Columns.Item(13).ColumnWidth = 10.666666666666666

$ws.Range("A6:XFD6").Select() | Out-Null
